$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new IPO record (삼성스팩9호) needs to be inserted as row 2 of the data
# table (sheet row 2, just under the header row). Every existing data row
# (originally rows 2-27) shifts down by one (to rows 3-28).
#
# The COM-interop Range/Rows .Insert() here stamps an inherited style index
# onto the new row's cells (it copies formatting from the row above), which
# the target file does not have (plain, unstyled data cells). So instead we
# shift the existing data down manually, cell by cell, from the bottom up,
# then populate row 2 with the new record's values. `.Value` is unreliable
# as a *getter* in this runtime (it returns a stub descriptor string rather
# than the underlying value) so reads use `.Value2`; writes use `.Value`.

$lastDataRow = 27
$lastCol = 20

for ($r = $lastDataRow; $r -ge 2; $r--) {
    $destRow = $r + 1
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $ws.Cells.Item($r, $c).Value2
    }
}

# Populate the newly freed row 2 with the 삼성스팩9호 record.
$ws.Cells.Item(2, 1).Value = "2023-11-23"
$ws.Cells.Item(2, 2).Value = "삼성스팩9호"
$ws.Cells.Item(2, 3).Value = "삼성"
$ws.Cells.Item(2, 4).Value = "2023-11-28"
$ws.Cells.Item(2, 5).Value = "2023-12-04"
$ws.Cells.Item(2, 6).Value = 20000000
$ws.Cells.Item(2, 7).Value = 10000000
$ws.Cells.Item(2, 8).Value = "-"
$ws.Cells.Item(2, 9).Value = 2000
$ws.Cells.Item(2, 10).Value = 2000
$ws.Cells.Item(2, 11).Value = "-"
$ws.Cells.Item(2, 12).Value = 2000
$ws.Cells.Item(2, 13).Value = "-"
$ws.Cells.Item(2, 14).Value = "-"
$ws.Cells.Item(2, 15).Value = 0
$ws.Cells.Item(2, 16).Value = "-"
$ws.Cells.Item(2, 17).Value = "-"
$ws.Cells.Item(2, 18).Value = "217.87 : 1"
$ws.Cells.Item(2, 19).Value = "-"
$ws.Cells.Item(2, 20).Value = "-"
